$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell value updates -------------------------------------------------
# Header F1: "remark" -> "Expected "
$ws.Range("F1").Value = "Expected "

# Row 4 test-data rename ("fake5" -> "fake10") + typo fix on the phone number
$ws.Range("A4").Value = "fake10"
$ws.Range("B4").Value = "fake10@g.com"
$ws.Range("C4").Value = 987654321

# --- Header row formatting (bold + yellow fill) -------------------------
# Build the combined format as a single named style so only ONE new font,
# ONE new fill and ONE new cellXf record are produced, then drop the name
# so the workbook keeps using a direct (unnamed) cell format - matching a
# plain "bold + yellow fill" direct formatting edit.
$headerStyle = $wb.Styles.Add("__HeaderFormatTemp")
$headerStyle.Font.Bold = $true
$headerStyle.Interior.Color = 65535
$ws.Range("A1:F1").Style = "__HeaderFormatTemp"
$wb.Styles.Item("__HeaderFormatTemp").Delete()

# --- Page setup: portrait orientation -----------------------------------
$ws.PageSetup.Orientation = 1

# --- Selection change (cosmetic, matches last-saved cursor position) ----
$ws.Range("E8").Select()
